# Updated 24V test cases and test data with new loading details method
# Adds a new "Loading Details Name" column (G) to the "Add Panels" sheet,
# populating the header and the per-row "Main Processor 24V (A)" label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell G7: "Loading Details Name" ---
# Clone formatting from the existing header cell F7 (bold font / blue fill / border)
# so the new header matches the look of the other headers, then set its text.
$ws.Range("F7").Copy($ws.Range("G7"))
$ws.Range("G7").Value = "Loading Details Name"

# --- New data cells G8:G13: "Main Processor 24V (A)" ---
# Clone formatting from an existing data cell in the same rows (A8 uses the
# grey fill + left alignment shared by the row) and additionally wrap text,
# matching the new style used for this column.
$ws.Range("A8").Copy($ws.Range("G8:G13"))
$ws.Range("G8:G13").WrapText = $true

$ws.Range("G8").Value = "Main Processor 24V (A)"
$ws.Range("G9").Value = "Main Processor 24V (A)"
$ws.Range("G10").Value = "Main Processor 24V (A)"
$ws.Range("G11").Value = "Main Processor 24V (A)"
$ws.Range("G12").Value = "Main Processor 24V (A)"
$ws.Range("G13").Value = "Main Processor 24V (A)"

# Size the new column to fit its contents, mirroring the other bestFit columns.
$ws.Columns("G:G").AutoFit()

# Leave the selection where the editor finished working, on G10.
$null = $ws.Range("G10").Select()
